# SectorGroup.xlsx: columns E ("codeforiati:category-name") and F
# ("codeforiati:group-code") were swapped — column E should actually hold
# the group-code and column F should hold the category-name (the header
# labels move along with the data). Swap the two columns' contents for
# every row, including the header row.
#
# The swap is done with Copy / PasteSpecial(xlPasteValues) instead of a
# plain Value assignment so that numeric-looking text codes such as
# "110"/"120" keep their original text (shared-string) cell type instead
# of Excel re-inferring them as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

$colE = "E" + $firstRow + ":E" + $lastRow
$colF = "F" + $firstRow + ":F" + $lastRow
$colTmp = "Z" + $firstRow + ":Z" + $lastRow

$xlPasteValues = -4163

# Stash column F's original values in a scratch column (Z).
$ws.Range($colF).Copy()
$ws.Range($colTmp).PasteSpecial($xlPasteValues)

# Move column E's values into column F.
$ws.Range($colE).Copy()
$ws.Range($colF).PasteSpecial($xlPasteValues)

# Move the stashed original column F values into column E.
$ws.Range($colTmp).Copy()
$ws.Range($colE).PasteSpecial($xlPasteValues)

# Clean up the scratch column and clipboard state.
$ws.Range($colTmp).Clear()
$excel.CutCopyMode = 0
